$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row (row 1) column titles
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Normalize capitalization of "de"/"el"/"la" -> "De"/"El"/"La" in a few cells
$ws.Range("B6").Value = "Mazapa De Madero"
$ws.Range("A9").Value = "Ciudad De México"
$ws.Range("A11").Value = "Estado De México"
$ws.Range("B13").Value = "Apaseo El Alto"
$ws.Range("B14").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B24").Value = "Teocuitatlán De Corona"

# 3) Remove trailing footer/metadata rows (45-49 and 476-480), which also
#    shrinks the sheet's used range/dimension down to A1:D43.
$ws.Range("A476:A480").EntireRow.Delete()
$ws.Range("A45:A49").EntireRow.Delete()
